$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Test1"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 45860
$ws.Range("B3").Value = "sedrftgyhuioygtfrd"
